$d = $word.ActiveDocument

# The first two paragraphs ("Feast of St. Joseph" / "By Dorothy Day") are
# replaced by a pandoc-style title block: a Title-styled heading whose
# words are individual runs, followed by an Authors-styled byline (also
# split into individual runs), with the old bookmark wrapping the heading
# removed entirely.

$secondPara = $d.Paragraphs.Item(2)
$headRange = $d.Range(0, $secondPara.Range.End)
$headRange.Delete()

# The original bookmark ("feast-of-st.-joseph") sat right at the very start
# of the document, around the heading paragraph. Its name isn't a legal
# Word bookmark identifier (it contains '.'), so it doesn't show up in the
# Bookmarks collection, but its start/end markers are still anchored at
# position 0. Deleting a zero-length range there twice removes the
# leftover bookmarkStart and bookmarkEnd markers.
$d.Range(0, 0).Delete()
$d.Range(0, 0).Delete()

$titleBlockXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t xml:space="preserve">Feast</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">of</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">St</w:t></w:r><w:r><w:t xml:space="preserve">.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Joseph</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Authors"/></w:pPr><w:r><w:t xml:space="preserve">Dorothy</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Day</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$d.Range(0, 0).InsertXML($titleBlockXml) | Out-Null
